# Applies the content edits described by the commit diff to the
# "Complaints Procedure" document:
#   - fixes an accidental double space
#   - updates the setting's phone number
#   - capitalises "ofsted" -> "Ofsted" (and lets the now-correct spelling
#     clear the old spell-check markers)
#   - refreshes the "Date" / "Review Date" lines at the bottom of the
#     document

$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindWrap (wrap) = 1
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($find, $replace, [bool]$matchCase = $false) {
    $d.Content.Find.Execute(
        $find, $matchCase, $false, $false, $false, $false,
        $true, $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

# Stray double space before "record will be stored in the Complaints Log."
Replace-Text "The  record will be stored" "The record will be stored"

# Updated OFSTED contact telephone number.
Replace-Text "Tel: 0300 123 1231" "Tel: 0300 123 4666"

# Lower-cased "ofsted" should read "Ofsted" (match case so the correctly
# capitalised occurrences elsewhere, and the all-caps "OFSTED" heading,
# are left untouched).
Replace-Text "ofsted" "Ofsted" $true

# Sign-off / review dates.
Replace-Text "Date: 31st August 2022" "Date: 14th September 2025"
Replace-Text "Review Date: August 2023" "Review Date: September 2026"
